# Append new Lancers job-listing rows and refresh the "取得日時" timestamp.
#
# Target final state for sheet "ランサーズ" (Worksheets.Item(1)):
#   - A2:A10 timestamp refreshed to 2025-11-05 12:38:29
#   - a new job (5427648) now occupies row 6, pushing the previous
#     rows 6-8 (5427397 / 5427338 / 5427459) down to rows 7-9
#   - a new job (5427699) is appended as row 10
#   - dimension grows from A1:H8 to A1:H10
#   - hyperlinks on column F are rebuilt so F2:F10 all point at the
#     correct URL for their row
#
# Implementation note: rather than calling Rows.Insert() (which, in this
# host, leaves the worksheet's Hyperlinks ref/target pairing stale - the
# text shifts down a row but the hyperlink keeps pointing at its old
# target), every data cell for rows 2-10 is written explicitly with its
# final value, and the F-column hyperlinks are deleted/recreated from
# scratch afterwards so every ref/target pair is correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-11-05 12:38:29"

# Final row data (row, title, category, price, deadline, url, score, skill-summary)
$data = @(
    @(2,  "医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5416301", 385, "🔥AI,Ai ◆開発 ◇アプリ"),
    @(3,  "ワードプレスサイト内に、chatgptのテキスト自動作成と自動でコピー状態の設定", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420440", 350, "🔥GPT,ChatGPT ◇サイト"),
    @(4,  "Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5416328", 310, "🔥AI,Ai"),
    @(5,  "<Next.js、バックエンド開発> ガントチャートアプリの改修製造", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427011", 225, "🔥Next.js ◆開発 ◇アプリ"),
    @(6,  "【急募】キントーン見積書をエクセルに変換するツール開発", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427648", 120, "◆ツール,開発"),
    @(7,  "デフォルトカメラ機能を活用したアプリ開発", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427397", 100, "◆開発 ◇アプリ"),
    @(8,  "弥生販売 得意先台帳登録 商品登録 売上伝票作成ツールのご相談", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427338", 73, "◆ツール"),
    @(9,  "【継続依頼あり】GASやn8nのオンラインセミナー研修講師を募集!", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427459", 13, ""),
    @(10, "中国語 ワードプレスの分かる人", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427699", 10, "")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne "") {
        $ws.Cells.Item($r, 8).Value = $row[7]
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
}

# Rebuild the column-F hyperlinks from scratch: deleting a single hyperlink
# in this host clears the whole-sheet collection as a side effect, so do it
# once up front and then re-add every row's link in the correct final order.
$ws.Range("F2").Hyperlinks.Delete()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $data[$i][0]
    $url = $data[$i][5]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url)
}
